# shark attack added + player made smaller to make shark scarier
#
# 1) Re-split the "Shark was rotating wrong ..." paragraph so the word
#    "player" is wrapped in its own run flanked by proofErr gramStart/gramEnd
#    markers (instead of the original split which bracketed "degrees.").
# 2) Add a blank paragraph after it.
# 3) Add a new paragraph about removing the destructible blocks.
# 4) Add a trailing blank paragraph.

$d = $word.ActiveDocument

$targetText = "Shark was rotating wrong when chasing player it was off by 90 degrees, fixed this by manually editing the picture 90 degrees."

$r = $d.Content
$found = $r.Find.Execute($targetText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the target 'Shark was rotating wrong ...' paragraph text"
}

# Clear the matched range; it collapses to an insertion point in its place.
$r.Text = ""

$newXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Shark was rotating wrong when chasing </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>player</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> it was off by 90 degrees, fixed this by manually editing the picture 90 degrees.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Made the decision to remove the destructible blocks as they blocked the shark and this way I could have the shark follow you from the start </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
'@

[void]$r.InsertXML($newXml)
